# Weekly update: insert a new price record for "Poroto verde" at
# Feria Lagunitas de Puerto Montt as the newest (first) entry in the
# data block, shifting the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 32 (the first data row of this
# sub-series); this shifts rows 32..51 down to 33..52 and copies the
# formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly observation.
$ws.Cells.Item(32, 1).Value  = 4
$ws.Cells.Item(32, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(32, 3).Value  = "Los Lagos"
$ws.Cells.Item(32, 4).Value  = 44523
$ws.Cells.Item(32, 5).Value  = 10
$ws.Cells.Item(32, 6).Value  = 100112031
$ws.Cells.Item(32, 7).Value  = "Poroto verde"
$ws.Cells.Item(32, 8).Value  = "Magnum"
$ws.Cells.Item(32, 9).Value  = "Primera"
$ws.Cells.Item(32, 10).Value = 35
$ws.Cells.Item(32, 11).Value = 43000
$ws.Cells.Item(32, 12).Value = 43000
$ws.Cells.Item(32, 13).Value = 43000
$ws.Cells.Item(32, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 1720
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
